$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.688.91"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "3.234.43"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.00%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "3.232.98"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("E12").Value = "  -3.23%  "
$ws.Range("D13").Value = "3.791.19"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("E14").Value = "  -3.11%  "
$ws.Range("D15").Value = "64.795.82"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "3.234.00"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "414.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("E26").Value = "  +4.76%  "
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000110"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  -5.10%  "
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").Value = "2.805.98"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  -6.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0630"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("E47").Value = "  -4.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "302.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
